$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.378.68"
$ws.Range("E2").Value = "  -4.60%  "

$ws.Range("D3").Value = "'3.358.80"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'568.11"
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").Value = "'132.32"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'3.359.66"
$ws.Range("E8").Value = "  -1.84%  "

$ws.Range("D9").Value = "'0.475"
$ws.Range("E9").Value = "  -1.06%  "

$ws.Range("D10").Value = "'7.52"
$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("E11").Value = "  -2.47%  "

$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Value = "'3.927.89"
$ws.Range("E13").Value = "  -1.90%  "

$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000171"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.350.46"
$ws.Range("E16").Value = "  -2.17%  "

$ws.Range("D17").Value = "'24.89"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "'60.481.32"
$ws.Range("E18").Value = "  -4.53%  "

$ws.Range("D19").Value = "'13.63"
$ws.Range("E19").Value = "  +3.36%  "

$ws.Range("E20").Value = "  -4.44%  "

$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").Value = "'365.30"
$ws.Range("E22").Value = "  -4.42%  "

$ws.Range("D23").Value = "'0.561"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'3.485.96"
$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "'69.77"
$ws.Range("E26").Value = "  -5.12%  "

$ws.Range("E27").Value = "  +3.73%  "

$ws.Range("E28").Value = "  +17.31%  "

$ws.Range("D29").Value = "'7.55"
$ws.Range("E29").Value = "  +7.47%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Value = "'0.154"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "'3.386.80"
$ws.Range("E35").Value = "  -1.87%  "

$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").Value = "'5.38"
$ws.Range("E37").Value = "  +3.76%  "

$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("D40").Value = "'158.98"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("E41").Value = "  +2.02%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "'1.21"
$ws.Range("E43").Value = "  +10.38%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.40"
$ws.Range("E44").Value = "  +2.46%  "

$ws.Range("D45").Value = "'40.98"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").Value = "'0.753"
$ws.Range("E46").Value = "  -3.85%  "

$ws.Range("D47").Value = "'23.89"
$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("D48").Value = "'1.59"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").Value = "'22.75"
$ws.Range("E50").Value = "  +12.36%  "

$ws.Range("D51").Value = "'0.895"
$ws.Range("E51").Value = "  +1.51%  "
